$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits (shared-string text updates) ---
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Row heights: 18.75 -> 19.5 for header + both data rows ---
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5

# --- Formatting: extend the bordered numeric style already used by
# I2/K2/I3/K3 up into the header cells I1/K1 (matches them onto the same
# style bucket), then make that font's color an explicit black instead of
# the implicit theme color. ---
$ws.Range("I1").Font.Name = "Calibri"
$ws.Range("K1").Font.Name = "Calibri"

$ws.Range("I1").Font.Color = 0
$ws.Range("K1").Font.Color = 0
$ws.Range("I2").Font.Color = 0
$ws.Range("K2").Font.Color = 0
$ws.Range("I3").Font.Color = 0
$ws.Range("K3").Font.Color = 0
